$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 425
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45180
